# Update "想去人数" (want-to-go count) values in column F across sheets
# as published by the latest gh-pages data refresh (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1503
$ws1.Range("F4").Value = 2113
$ws1.Range("F5").Value = 7607
$ws1.Range("F6").Value = 4828
$ws1.Range("F7").Value = 7108
$ws1.Range("F10").Value = 1506
$ws1.Range("F11").Value = 868
$ws1.Range("F13").Value = 60
$ws1.Range("F14").Value = 1171
$ws1.Range("F16").Value = 177
$ws1.Range("F20").Value = 234
$ws1.Range("F22").Value = 1190
$ws1.Range("F23").Value = 957
$ws1.Range("F24").Value = 4
$ws1.Range("F31").Value = 195
$ws1.Range("F33").Value = 2
$ws1.Range("F34").Value = 49
$ws1.Range("F35").Value = 115
$ws1.Range("F36").Value = 37
$ws1.Range("F37").Value = 554
$ws1.Range("F38").Value = 428
$ws1.Range("F39").Value = 78
$ws1.Range("F41").Value = 89
$ws1.Range("F42").Value = 395
$ws1.Range("F43").Value = 1203
$ws1.Range("F45").Value = 146

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 133
$ws2.Range("F17").Value = 560
$ws2.Range("F22").Value = 211
$ws2.Range("F32").Value = 862
$ws2.Range("F33").Value = 993
$ws2.Range("F34").Value = 612
$ws2.Range("F41").Value = 144

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value = 48
$ws3.Range("F9").Value = 52
$ws3.Range("F10").Value = 1653
$ws3.Range("F11").Value = 2552

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1503
$ws4.Range("F9").Value = 7607
$ws4.Range("F10").Value = 4828
$ws4.Range("F11").Value = 7108
$ws4.Range("F13").Value = 1506
$ws4.Range("F15").Value = 868
$ws4.Range("F18").Value = 1653
$ws4.Range("F19").Value = 2552
$ws4.Range("F20").Value = 211
$ws4.Range("F21").Value = 60
$ws4.Range("F22").Value = 1171
$ws4.Range("F23").Value = 177
$ws4.Range("F25").Value = 234
$ws4.Range("F26").Value = 1190
$ws4.Range("F28").Value = 957
$ws4.Range("F31").Value = 195
$ws4.Range("F33").Value = 863
$ws4.Range("F34").Value = 49
$ws4.Range("F35").Value = 115
$ws4.Range("F36").Value = 993
$ws4.Range("F37").Value = 554
$ws4.Range("F38").Value = 612
$ws4.Range("F39").Value = 78
$ws4.Range("F41").Value = 89
$ws4.Range("F43").Value = 395
$ws4.Range("F47").Value = 146
